$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$survey.Range("B2").Value = "select_one countries"
$survey.Range("F2").Value = "country"
$survey.Range("G2").Value = "Choose a country:"

$survey.Range("B3").Value = "select_one states"
$survey.Range("F3").Value = "state"
$survey.Range("G3").Value = "Choose a state:"

$queries = $wb.Worksheets.Item("queries")
$queries.Range("A2").Value = "countries"

$queries.Range("A3").Value = "states"
$queries.Range("B3").Value = "`"https://query.yahooapis.com/v1/public/yql?format=json&q=`" +  encodeURIComponent(`"select * from geo.states where place='`" + data('country') + `"'`")"
